$wb = $excel.ActiveWorkbook

$wsAssets = $wb.Worksheets.Item("Assets")
$wsLiabilities = $wb.Worksheets.Item("Liabilities")

# --- Assets sheet (sheet1): update transaction classification totals ---
$wsAssets.Range("B2").Value = 57327.05
$wsAssets.Range("C2").Value = 69089.91
$wsAssets.Range("D2").Value = 46505.97
$wsAssets.Range("E2").Value = 68042.99000000001
$wsAssets.Range("F2").Value = 97647.47
$wsAssets.Range("G2").Value = 338613.39
$wsAssets.Range("B3").Value = 452802.19
$wsAssets.Range("C3").Value = 439948.04
$wsAssets.Range("D3").Value = 521104.54
$wsAssets.Range("E3").Value = 400785.61
$wsAssets.Range("F3").Value = 439761.29
$wsAssets.Range("G3").Value = 2254401.67
$wsAssets.Range("B4").Value = -11891.18
$wsAssets.Range("C4").Value = -20277.56
$wsAssets.Range("D4").Value = -16132.08
$wsAssets.Range("E4").Value = -18620.19
$wsAssets.Range("F4").Value = -13311.51
$wsAssets.Range("G4").Value = -80232.52
$wsAssets.Range("B5").Value = -28564.63
$wsAssets.Range("C5").Value = -29994.23
$wsAssets.Range("D5").Value = -27644.43
$wsAssets.Range("E5").Value = -29930.08
$wsAssets.Range("F5").Value = -24965.85
$wsAssets.Range("G5").Value = -141099.22
$wsAssets.Range("B6").Value = 9349.23
$wsAssets.Range("C6").Value = 9672.219999999999
$wsAssets.Range("D6").Value = 7688.38
$wsAssets.Range("E6").Value = 8537.84
$wsAssets.Range("F6").Value = 19512.25
$wsAssets.Range("G6").Value = 54759.92
$wsAssets.Range("B7").Value = -196582.08
$wsAssets.Range("C7").Value = -188143.44
$wsAssets.Range("D7").Value = -185434.74
$wsAssets.Range("E7").Value = -205682.74
$wsAssets.Range("F7").Value = -189920.85
$wsAssets.Range("G7").Value = -965763.85
$wsAssets.Range("B8").Value = -18135.67
$wsAssets.Range("C8").Value = -16482.58
$wsAssets.Range("D8").Value = -15629.56
$wsAssets.Range("E8").Value = -19907.67
$wsAssets.Range("F8").Value = -12846.96
$wsAssets.Range("G8").Value = -83002.44
$wsAssets.Range("B9").Value = 81716.10000000001
$wsAssets.Range("C9").Value = 76255.21000000001
$wsAssets.Range("D9").Value = 76187.89999999999
$wsAssets.Range("E9").Value = 68580.99000000001
$wsAssets.Range("F9").Value = 110584.34
$wsAssets.Range("G9").Value = 413324.54
$wsAssets.Range("B10").Value = 346021.01
$wsAssets.Range("C10").Value = 340067.57
$wsAssets.Range("D10").Value = 406645.98
$wsAssets.Range("E10").Value = 271806.75
$wsAssets.Range("F10").Value = 426460.18
$wsAssets.Range("G10").Value = 1791001.49

# --- Liabilities sheet (sheet2): update transaction classification totals ---
$wsLiabilities.Range("B2").Value = -38915.87
$wsLiabilities.Range("C2").Value = -42885.77
$wsLiabilities.Range("D2").Value = -38293.06
$wsLiabilities.Range("E2").Value = -31954.07
$wsLiabilities.Range("F2").Value = -28973.95
$wsLiabilities.Range("G2").Value = -181022.72
$wsLiabilities.Range("B3").Value = -36353.05
$wsLiabilities.Range("C3").Value = -30457.38
$wsLiabilities.Range("D3").Value = -38819.41
$wsLiabilities.Range("E3").Value = -38946.09
$wsLiabilities.Range("F3").Value = -34430.07
$wsLiabilities.Range("G3").Value = -179006
$wsLiabilities.Range("B4").Value = -11866.61
$wsLiabilities.Range("C4").Value = -20344.57
$wsLiabilities.Range("D4").Value = -15438.74
$wsLiabilities.Range("E4").Value = -15008.15
$wsLiabilities.Range("F4").Value = -18061.25
$wsLiabilities.Range("G4").Value = -80719.32000000001
$wsLiabilities.Range("B5").Value = -42160.42
$wsLiabilities.Range("C5").Value = -40308.66
$wsLiabilities.Range("D5").Value = -31167.05
$wsLiabilities.Range("E5").Value = -37792.51
$wsLiabilities.Range("F5").Value = -39080.4
$wsLiabilities.Range("G5").Value = -190509.04
$wsLiabilities.Range("B6").Value = -57936.9
$wsLiabilities.Range("C6").Value = -44474.84
$wsLiabilities.Range("D6").Value = -54028.1
$wsLiabilities.Range("E6").Value = -73444.08
$wsLiabilities.Range("F6").Value = -49862.93
$wsLiabilities.Range("G6").Value = -279746.85
$wsLiabilities.Range("B7").Value = 187232.85
$wsLiabilities.Range("C7").Value = 178471.22
$wsLiabilities.Range("D7").Value = 177746.36
$wsLiabilities.Range("E7").Value = 197144.9
$wsLiabilities.Range("F7").Value = 170408.6
$wsLiabilities.Range("G7").Value = 911003.9300000001
